$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing header style (already applied to B1 / A2) as the template
# for the newly-added header cells, via copy / paste-special-formats.
$headerTemplate = $ws.Range("B1")

# --- New "from" header columns: CHP1, pvt1 ---
$ws.Range("E1").Value = "P_from_CHP1"
$headerTemplate.Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("F1").Value = "P_from_pvt1"
$headerTemplate.Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- New "to" header row: charging_station2 ---
$ws.Range("A6").Value = "P_to_charging_station2"
$headerTemplate.Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B6").Value = "P_net1_charging_station2"
$ws.Range("C6").Value = "P_pv1_charging_station2"
$ws.Range("D6").Value = "P_bat1_charging_station2"

# --- New column E body: CHP1 -> * ---
$ws.Range("E2").Value = "P_CHP1_demand1"
$ws.Range("E3").Value = "P_CHP1_net1"
$ws.Range("E4").Value = "P_CHP1_bat1"
$ws.Range("E5").Value = "P_CHP1_charging_station1"
$ws.Range("E6").Value = "P_CHP1_charging_station2"

# --- New column F body: pvt1 -> * ---
$ws.Range("F2").Value = "P_pvt1_demand1"
$ws.Range("F3").Value = "P_pvt1_net1"
$ws.Range("F4").Value = "P_pvt1_bat1"
$ws.Range("F5").Value = "P_pvt1_charging_station1"
$ws.Range("F6").Value = "P_pvt1_charging_station2"

$excel.CutCopyMode = $false
